$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts Player ID, Player, etc. one column right)
$ws.Columns("A").Insert()

# Fill in the new "Match ID" column header (row 2, the real header row)
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Row 3 is a hidden spacer row; give its A cell the same bold style without a value.
# Toggle visibility around the write so the engine doesn't stamp a custom row height
# on a hidden row.
$ws.Rows(3).Hidden = $false
$ws.Range("A3").Font.Bold = $true
$ws.Rows(3).Hidden = $true

# Data rows 4-19 all belong to match id 9
$ws.Range("A4:A19").Value = 9
$ws.Range("A4:A19").Font.Bold = $true

# Row 20 is the hidden totals row; same visibility toggle trick, no bold styling.
$ws.Rows(20).Hidden = $false
$ws.Range("A20").Value = 9
$ws.Rows(20).Hidden = $true

# Update the selection to match the committed worksheet state
$ws.Range("A2:A19").Select()
